$wb = $excel.ActiveWorkbook

# This script updates specific cell values in the profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as refreshed by the scheduled
# market-data runner job.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 278014.3
$ws.Range("I137").Value = 384975.56
$ws.Range("J137").Value = 1697.75
$ws.Range("K137").Value = 1154926.68
$ws.Range("L137").Value = 5093.25
$ws.Range("M137").Value = -1152376.68
$ws.Range("N137").Value = -10193.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7306.2354
$ws.Range("I61").Value = 9103
$ws.Range("J61").Value = 2994
$ws.Range("K61").Value = 9103
$ws.Range("L61").Value = 2994
$ws.Range("M61").Value = -8891
$ws.Range("N61").Value = -3418
$ws.Range("H74").Value = 6660.5264
$ws.Range("I74").Value = 1088.2222
$ws.Range("J74").Value = 11675.6
$ws.Range("K74").Value = 1088.2222
$ws.Range("L74").Value = 11675.6
$ws.Range("M74").Value = -214.2221999999999
$ws.Range("N74").Value = -13423.6
$ws.Range("H77").Value = 6660.5264
$ws.Range("I77").Value = 1088.2222
$ws.Range("J77").Value = 11675.6
$ws.Range("K77").Value = 5441.111
$ws.Range("L77").Value = 58378
$ws.Range("M77").Value = -1073.111
$ws.Range("N77").Value = -67114
$ws.Range("H122").Value = 42858056
$ws.Range("I122").Value = 50000956
$ws.Range("J122").Value = 653.5
$ws.Range("K122").Value = 150002868
$ws.Range("L122").Value = 1960.5
$ws.Range("M122").Value = -150000418
$ws.Range("N122").Value = -6860.5
$ws.Range("H132").Value = 8336158
$ws.Range("I132").Value = 15626984
$ws.Range("J132").Value = 3785.2856
$ws.Range("K132").Value = 46880952
$ws.Range("L132").Value = 11355.8568
$ws.Range("M132").Value = -46878422
$ws.Range("N132").Value = -16415.8568
$ws.Range("H136").Value = 7306.2354
$ws.Range("I136").Value = 9103
$ws.Range("J136").Value = 2994
$ws.Range("K136").Value = 27309
$ws.Range("L136").Value = 8982
$ws.Range("M136").Value = -24759
$ws.Range("N136").Value = -14082

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 68853.336
$ws.Range("J59").Value = 68853.336
$ws.Range("L59").Value = 68853.336
$ws.Range("N59").Value = -70547.336
$ws.Range("H74").Value = 27495
$ws.Range("J74").Value = 27495
$ws.Range("L74").Value = 27495
$ws.Range("N74").Value = -29367
$ws.Range("H77").Value = 27495
$ws.Range("J77").Value = 27495
$ws.Range("L77").Value = 82485
$ws.Range("N77").Value = -91845
$ws.Range("H134").Value = 33386436
$ws.Range("I134").Value = 37095484
$ws.Range("K134").Value = 111286452
$ws.Range("M134").Value = -111283917

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9617.125
$ws.Range("I31").Value = 886.6061
$ws.Range("J31").Value = 28824.268
$ws.Range("K31").Value = 886.6061
$ws.Range("L31").Value = 28824.268
$ws.Range("M31").Value = -591.6061
$ws.Range("N31").Value = -29414.268
$ws.Range("H34").Value = 9617.125
$ws.Range("I34").Value = 886.6061
$ws.Range("J34").Value = 28824.268
$ws.Range("K34").Value = 886.6061
$ws.Range("L34").Value = 28824.268
$ws.Range("M34").Value = -684.6061
$ws.Range("N34").Value = -29228.268
$ws.Range("H58").Value = 4366325.5
$ws.Range("I58").Value = 7193706.5
$ws.Range("J58").Value = 16508.309
$ws.Range("K58").Value = 7193706.5
$ws.Range("L58").Value = 16508.309
$ws.Range("M58").Value = -7193503.5
$ws.Range("N58").Value = -16914.309
$ws.Range("H132").Value = 9528741
$ws.Range("I132").Value = 18519938
$ws.Range("J132").Value = 8648.588
$ws.Range("K132").Value = 55559814
$ws.Range("L132").Value = 25945.764
$ws.Range("M132").Value = -55557284
$ws.Range("N132").Value = -31005.764
$ws.Range("H134").Value = 11162348
$ws.Range("I134").Value = 11365179
$ws.Range("J134").Value = 10418632
$ws.Range("K134").Value = 34095537
$ws.Range("L134").Value = 31255896
$ws.Range("M134").Value = -34093002
$ws.Range("N134").Value = -31260966
$ws.Range("H136").Value = 4366325.5
$ws.Range("I136").Value = 7193706.5
$ws.Range("J136").Value = 16508.309
$ws.Range("K136").Value = 21581119.5
$ws.Range("L136").Value = 49524.927
$ws.Range("M136").Value = -21578569.5
$ws.Range("N136").Value = -54624.927

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3794.6
$ws.Range("J80").Value = 4368.125
$ws.Range("L80").Value = 13104.375
$ws.Range("N80").Value = -14976.375
$ws.Range("H83").Value = 3794.6
$ws.Range("J83").Value = 4368.125
$ws.Range("L83").Value = 39313.125
$ws.Range("N83").Value = -48673.125
$ws.Range("H113").Value = 4112.5713
$ws.Range("I113").Value = 541.25
$ws.Range("J113").Value = 8874.333000000001
$ws.Range("K113").Value = 1623.75
$ws.Range("L113").Value = 26622.999
$ws.Range("M113").Value = 546.25
$ws.Range("N113").Value = -30962.999
$ws.Range("H122").Value = 578.4761999999999
$ws.Range("I122").Value = 437.88235
$ws.Range("J122").Value = 1176
$ws.Range("K122").Value = 3940.94115
$ws.Range("L122").Value = 10584
$ws.Range("M122").Value = -1490.94115
$ws.Range("N122").Value = -15484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 45455710
$ws.Range("I122").Value = 71429490
$ws.Range("J122").Value = 1598.5
$ws.Range("K122").Value = 214288470
$ws.Range("L122").Value = 4795.5
$ws.Range("M122").Value = -214286020
$ws.Range("N122").Value = -9695.5
$ws.Range("H132").Value = 142869360
$ws.Range("I132").Value = 500001000
$ws.Range("J132").Value = 16700.6
$ws.Range("K132").Value = 1500003000
$ws.Range("L132").Value = 50101.8
$ws.Range("M132").Value = -1500000470
$ws.Range("N132").Value = -55161.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 59527164
$ws.Range("I122").Value = 83335570
$ws.Range("K122").Value = 250006710
$ws.Range("M122").Value = -250004260
$ws.Range("H132").Value = 7694446.5
$ws.Range("I132").Value = 11113034
$ws.Range("J132").Value = 2624.875
$ws.Range("K132").Value = 33339102
$ws.Range("L132").Value = 7874.625
$ws.Range("M132").Value = -33336572
$ws.Range("N132").Value = -12934.625
$ws.Range("H136").Value = 4156.0234
$ws.Range("I136").Value = 4772.697
$ws.Range("J136").Value = 2121
$ws.Range("K136").Value = 14318.091
$ws.Range("L136").Value = 6363
$ws.Range("M136").Value = -11768.091
$ws.Range("N136").Value = -11463
$ws.Range("H137").Value = 36800
$ws.Range("J137").Value = 37000
$ws.Range("L137").Value = 37000
$ws.Range("N137").Value = -47200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 50413400
$ws.Range("I132").Value = 54548296
$ws.Range("J132").Value = 42832764
$ws.Range("K132").Value = 163644888
$ws.Range("L132").Value = 128498292
$ws.Range("M132").Value = -163642358
$ws.Range("N132").Value = -128503352
$ws.Range("H136").Value = 20905680
$ws.Range("I136").Value = 13226087
$ws.Range("J136").Value = 35716324
$ws.Range("K136").Value = 39678261
$ws.Range("L136").Value = 107148972
$ws.Range("M136").Value = -39675711
$ws.Range("N136").Value = -107154072
